# Scheduled data refresh: update market-price / profit columns across all Leve sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2857.4
$ws.Range("I28").Value = 2765
$ws.Range("J28").Value = 2996
$ws.Range("K28").Value = 2765
$ws.Range("L28").Value = 2996
$ws.Range("M28").Value = -2280
$ws.Range("N28").Value = -3966

$ws.Range("H86").Value = 90957720
$ws.Range("I86").Value = 125065110
$ws.Range("K86").Value = 125065110
$ws.Range("M86").Value = -125063987

$ws.Range("H89").Value = 90957720
$ws.Range("I89").Value = 125065110
$ws.Range("K89").Value = 625325550
$ws.Range("M89").Value = -625319934

$ws.Range("H100").Value = 2668
$ws.Range("I100").Value = 2668
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 2668
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -2127
$ws.Range("N100").ClearContents()

$ws.Range("H112").Value = 4412.9565
$ws.Range("I112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("M112").ClearContents()

$ws.Range("H132").Value = 194826.39
$ws.Range("I132").Value = 236022.08
$ws.Range("J132").Value = 20889
$ws.Range("K132").Value = 708066.24
$ws.Range("L132").Value = 62667
$ws.Range("M132").Value = -705536.24
$ws.Range("N132").Value = -67727

$ws.Range("H137").Value = 12823030
$ws.Range("I137").Value = 1896.8462
$ws.Range("K137").Value = 5690.5386
$ws.Range("M137").Value = -3140.5386

$ws.Range("H138").Value = 2303.7476
$ws.Range("J138").Value = 2458.738
$ws.Range("L138").Value = 7376.214
$ws.Range("N138").Value = -17656.214

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4231
$ws.Range("I45").Value = 3771.818
$ws.Range("K45").Value = 3771.818
$ws.Range("M45").Value = -3394.818

$ws.Range("H61").Value = 3343.6333
$ws.Range("I61").Value = 2976.2444
$ws.Range("K61").Value = 2976.2444
$ws.Range("M61").Value = -2764.2444

$ws.Range("H74").Value = 13890918
$ws.Range("I74").Value = 22729004
$ws.Range("J74").Value = 2498.1428
$ws.Range("K74").Value = 22729004
$ws.Range("L74").Value = 2498.1428
$ws.Range("M74").Value = -22728130
$ws.Range("N74").Value = -4246.1428

$ws.Range("H77").Value = 13890918
$ws.Range("I77").Value = 22729004
$ws.Range("J77").Value = 2498.1428
$ws.Range("K77").Value = 113645020
$ws.Range("L77").Value = 12490.714
$ws.Range("M77").Value = -113640652
$ws.Range("N77").Value = -21226.714

$ws.Range("H97").Value = 398.65
$ws.Range("I97").Value = 404.29413
$ws.Range("K97").Value = 404.29413
$ws.Range("M97").Value = 91.70587

$ws.Range("H122").Value = 3231.1162
$ws.Range("I122").Value = 2221.5
$ws.Range("J122").Value = 5561
$ws.Range("K122").Value = 6664.5
$ws.Range("L122").Value = 16683
$ws.Range("M122").Value = -4214.5
$ws.Range("N122").Value = -21583

$ws.Range("H132").Value = 14003.135
$ws.Range("I132").Value = 17785.695
$ws.Range("K132").Value = 53357.085
$ws.Range("M132").Value = -50827.085

$ws.Range("H136").Value = 3343.6333
$ws.Range("I136").Value = 2976.2444
$ws.Range("K136").Value = 8928.733200000001
$ws.Range("M136").Value = -6378.733200000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3448.4856
$ws.Range("I20").Value = 3181.4119
$ws.Range("J20").Value = 3700.7222
$ws.Range("K20").Value = 3181.4119
$ws.Range("L20").Value = 3700.7222
$ws.Range("M20").Value = -2934.4119
$ws.Range("N20").Value = -4194.7222

$ws.Range("H99").Value = 3047.6924
$ws.Range("I99").Value = 3202.05
$ws.Range("K99").Value = 3202.05
$ws.Range("M99").Value = -1704.05

$ws.Range("H105").Value = 3212.9375
$ws.Range("I105").Value = 3094.0715
$ws.Range("J105").Value = 4045
$ws.Range("K105").Value = 3094.0715
$ws.Range("L105").Value = 4045
$ws.Range("M105").Value = -1347.0715
$ws.Range("N105").Value = -7539

$ws.Range("H132").Value = 98482.75
$ws.Range("J132").Value = 98482.75
$ws.Range("L132").Value = 98482.75
$ws.Range("N132").Value = -108602.75

$ws.Range("H134").Value = 4063.1292
$ws.Range("I134").Value = 4031.9
$ws.Range("K134").Value = 12095.7
$ws.Range("M134").Value = -9560.700000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 20411420
$ws.Range("I31").Value = 25643288
$ws.Range("J31").Value = 7130
$ws.Range("K31").Value = 25643288
$ws.Range("L31").Value = 7130
$ws.Range("M31").Value = -25642993
$ws.Range("N31").Value = -7720

$ws.Range("H34").Value = 20411420
$ws.Range("I34").Value = 25643288
$ws.Range("J34").Value = 7130
$ws.Range("K34").Value = 25643288
$ws.Range("L34").Value = 7130
$ws.Range("M34").Value = -25643086
$ws.Range("N34").Value = -7534

$ws.Range("H52").Value = 44623.332
$ws.Range("J52").Value = 40000
$ws.Range("L52").Value = 40000
$ws.Range("N52").Value = -40588

$ws.Range("H105").Value = 695
$ws.Range("I105").Value = 403.75
$ws.Range("K105").Value = 403.75
$ws.Range("M105").Value = 1343.25

$ws.Range("H141").Value = 163453.42
$ws.Range("J141").Value = 168741.25
$ws.Range("L141").Value = 168741.25
$ws.Range("N141").Value = -179101.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 597.86664
$ws.Range("I6").Value = 355.07144
$ws.Range("J6").Value = 3997
$ws.Range("K6").Value = 1065.21432
$ws.Range("L6").Value = 11991
$ws.Range("M6").Value = -952.21432
$ws.Range("N6").Value = -12217

$ws.Range("H68").Value = 1944.75
$ws.Range("I68").Value = 1889.5
$ws.Range("J68").Value = 2000
$ws.Range("K68").Value = 5668.5
$ws.Range("L68").Value = 6000
$ws.Range("M68").Value = -4857.5
$ws.Range("N68").Value = -7622

$ws.Range("H71").Value = 1944.75
$ws.Range("I71").Value = 1889.5
$ws.Range("J71").Value = 2000
$ws.Range("K71").Value = 17005.5
$ws.Range("L71").Value = 18000
$ws.Range("M71").Value = -12949.5
$ws.Range("N71").Value = -26112

$ws.Range("H113").Value = 858.8214
$ws.Range("J113").Value = 854.1667
$ws.Range("L113").Value = 2562.5001
$ws.Range("N113").Value = -6902.5001

$ws.Range("H114").Value = 6704
$ws.Range("I114").Value = 1456
$ws.Range("J114").Value = 7360
$ws.Range("K114").Value = 4368
$ws.Range("L114").Value = 22080
$ws.Range("M114").Value = -1114
$ws.Range("N114").Value = -28588

$ws.Range("H134").Value = 9689.85
$ws.Range("I134").Value = 4138.231
$ws.Range("K134").Value = 12414.693
$ws.Range("M134").Value = -7344.692999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 30493.5

$ws.Range("H88").Value = 93000
$ws.Range("J88").Value = 93000
$ws.Range("L88").Value = 93000
$ws.Range("N88").Value = -93902

$ws.Range("H91").Value = 93000
$ws.Range("J91").Value = 93000
$ws.Range("L91").Value = 93000
$ws.Range("N91").Value = -96120

$ws.Range("H97").Value = 1917.7142
$ws.Range("I97").Value = 1404.25
$ws.Range("J97").Value = 4998.5
$ws.Range("K97").Value = 1404.25
$ws.Range("L97").Value = 4998.5
$ws.Range("M97").Value = -908.25
$ws.Range("N97").Value = -5990.5

$ws.Range("H132").Value = 3577.4666
$ws.Range("I132").Value = 2862.7058
$ws.Range("J132").Value = 4512.154
$ws.Range("K132").Value = 8588.117400000001
$ws.Range("L132").Value = 13536.462
$ws.Range("M132").Value = -6058.117400000001
$ws.Range("N132").Value = -18596.462

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3620
$ws.Range("I100").Value = 3384.2
$ws.Range("K100").Value = 3384.2
$ws.Range("M100").Value = -2843.2

$ws.Range("H104").Value = 73332.664
$ws.Range("J104").Value = 73332.664
$ws.Range("L104").Value = 73332.664
$ws.Range("N104").Value = -80320.664

$ws.Range("H122").Value = 33341170
$ws.Range("I122").Value = 58827830
$ws.Range("J122").Value = 12457.77
$ws.Range("K122").Value = 176483490
$ws.Range("L122").Value = 37373.31
$ws.Range("M122").Value = -176481040
$ws.Range("N122").Value = -42273.31

$ws.Range("H135").Value = 89817.55
$ws.Range("J135").Value = 89817.55
$ws.Range("L135").Value = 89817.55
$ws.Range("N135").Value = -99957.55

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I54").Value = 40000
$ws.Range("J54").Value = 29333
$ws.Range("K54").Value = 40000
$ws.Range("L54").Value = 29333
$ws.Range("M54").Value = -39480
$ws.Range("N54").Value = -30373

$ws.Range("H107").Value = 632.13635
$ws.Range("I107").Value = 406.64706
$ws.Range("K107").Value = 1219.94118
$ws.Range("M107").Value = 700.05882

$ws.Range("H136").Value = 3771.1428
$ws.Range("I136").Value = 2544.1667
$ws.Range("J136").Value = 5611.607
$ws.Range("K136").Value = 7632.500100000001
$ws.Range("L136").Value = 16834.821
$ws.Range("M136").Value = -5082.500100000001
$ws.Range("N136").Value = -21934.821
